# Update FuelPrices at 2025-04-14 02:42
#
# - Cell C25 switches from the date-only format (style index 3,
#   "YYYY-MM-DD") to the date/time format (style index 2,
#   "YYYY-MM-DD HH:MM:SS"); its value (45756) stays the same.
# - A new data row 26 is appended with a new day's prices, using the
#   date-only format (style index 3) that C25 used to have.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C25: keep its value, just change the number format (style) from
# date-only to date+time.
$ws.Range("C25").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 26 with the next day's fuel price data.
$ws.Range("A26").Value = 770.419
$ws.Range("B26").Value = 697.8920000000001
$ws.Range("C26").Value = 45757
$ws.Range("C26").NumberFormat = "YYYY-MM-DD"
